# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "Valor Mora" (F column) figures for the first ("2011") and last
# ("2105") period rows of the account-statement table are swapped:
# the oldest period now carries the value that used to belong to the
# most recent period, and vice versa.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 -> period 2011 (first data row): was 56000, now 44800
$ws.Range("F16").Value = 44800

# Row 22 -> period 2105 (last data row): was 44800, now 56000
$ws.Range("F22").Value = 56000
